# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 17 and Row 18 ("Periodo Mora" / "Valor Mora" columns) are updated:
#   - E17: "2108" -> "2107"
#   - E18: "2107" -> "2108"
#   - G17: 908526 -> 877803
#   - G18: 908526 -> 877803

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between the two rows.
$ws.Range("E17").Value = "2107"
$ws.Range("E18").Value = "2108"

# Update the "Valor Mora" amounts for both rows.
$ws.Range("G17").Value = 877803
$ws.Range("G18").Value = 877803
